# [Kadastro App] Kayıt silindi: 11133517
#
# The record whose "Kayıt No" (column A) equals 11133517 is removed from
# the workbook entirely. It shows up twice: once in the master "Kayitlar"
# roll-up sheet and once more in the per-district "Merkez İlçe" sheet
# (every record also gets mirrored onto its district's own sheet). Deleting
# the whole row - rather than just blanking it - shifts every following
# row up by one, which is why the sheets' used range / dimension shrinks
# by exactly one row (A1:F728 -> A1:F727 on "Kayitlar", A1:F186 -> A1:F185
# on "Merkez İlçe").

$wb = $excel.ActiveWorkbook

$recordId = "11133517"

$sheetNames = @("Kayitlar", "Merkez İlçe")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # "Kayıt No" lives in column A. Use Find to jump straight to a
    # candidate row, then confirm the cell's text is an exact match
    # (Find()'s default match mode is substring-based) before deleting.
    $targetRow = 0
    $candidate = $ws.Columns.Item(1).Find($recordId)
    if (($candidate -ne $null) -and ($ws.Cells.Item($candidate.Row, 1).Text -eq $recordId)) {
        $targetRow = $candidate.Row
    } else {
        # Fallback: linear scan, in case Find ever returns a partial hit.
        $lastRow = $ws.UsedRange.Rows.Count
        for ($r = 1; $r -le $lastRow; $r++) {
            if ($ws.Cells.Item($r, 1).Text -eq $recordId) {
                $targetRow = $r
                break
            }
        }
    }

    if ($targetRow -gt 0) {
        # Delete the entire row so every subsequent row shifts up by one.
        $ws.Rows.Item($targetRow).Delete()
    }
}
